$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1504.2609
$ws.Range("J129").Value = 1654.9
$ws.Range("L129").Value = 4964.700000000001
$ws.Range("N129").Value = -14964.7

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1835.76
$ws.Range("I2").Value = 749.7273
$ws.Range("K2").Value = 749.7273
$ws.Range("M2").Value = -636.7273
$ws.Range("H55").Value = 19302
$ws.Range("J55").Value = 19302
$ws.Range("L55").Value = 19302
$ws.Range("N55").Value = -19932
$ws.Range("H80").Value = 31652
$ws.Range("J80").Value = 31652
$ws.Range("L80").Value = 31652
$ws.Range("N80").Value = -33648
$ws.Range("H83").Value = 31652
$ws.Range("J83").Value = 31652
$ws.Range("L83").Value = 94956
$ws.Range("N83").Value = -104940
$ws.Range("H116").Value = 1835.76
$ws.Range("I116").Value = 749.7273
$ws.Range("K116").Value = 749.7273
$ws.Range("M116").Value = 1544.2727
$ws.Range("H132").Value = 1897.1538
$ws.Range("I132").Value = 1792.7435
$ws.Range("K132").Value = 5378.2305
$ws.Range("M132").Value = -2848.2305

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1835.76
$ws.Range("I3").Value = 749.7273
$ws.Range("K3").Value = 749.7273
$ws.Range("M3").Value = -635.7273
$ws.Range("H22").Value = 1113.5
$ws.Range("I22").Value = 1113.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1113.5
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -940.5
$ws.Range("H35").Value = 27203.5
$ws.Range("J35").Value = 27203.5
$ws.Range("L35").Value = 27203.5
$ws.Range("N35").Value = -27823.5
$ws.Range("H86").Value = 2474.2307
$ws.Range("I86").Value = 2881.1428
$ws.Range("J86").Value = 1999.5
$ws.Range("K86").Value = 2881.1428
$ws.Range("L86").Value = 1999.5
$ws.Range("M86").Value = -1758.1428
$ws.Range("N86").Value = -4245.5
$ws.Range("H89").Value = 2474.2307
$ws.Range("I89").Value = 2881.1428
$ws.Range("J89").Value = 1999.5
$ws.Range("K89").Value = 14405.714
$ws.Range("L89").Value = 9997.5
$ws.Range("M89").Value = -8789.714
$ws.Range("N89").Value = -21229.5
$ws.Range("H135").Value = 56724
$ws.Range("J135").Value = 56724
$ws.Range("L135").Value = 56724
$ws.Range("N135").Value = -66864

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 27476.666
$ws.Range("J41").Value = 27476.666
$ws.Range("L41").Value = 27476.666
$ws.Range("N41").Value = -28332.666
$ws.Range("H50").Value = 24096
$ws.Range("J50").Value = 24096
$ws.Range("L50").Value = 24096
$ws.Range("N50").Value = -25346
$ws.Range("H51").Value = 27951.5
$ws.Range("J51").Value = 27951.5
$ws.Range("L51").Value = 27951.5
$ws.Range("N51").Value = -29423.5
$ws.Range("H59").Value = 32310
$ws.Range("J59").Value = 32310
$ws.Range("L59").Value = 32310
$ws.Range("N59").Value = -34600
$ws.Range("H60").Value = 20934.334
$ws.Range("J60").Value = 28401.5
$ws.Range("L60").Value = 28401.5
$ws.Range("N60").Value = -29423.5
$ws.Range("H61").Value = 27951.5
$ws.Range("J61").Value = 27951.5
$ws.Range("L61").Value = 27951.5
$ws.Range("N61").Value = -28647.5
$ws.Range("H62").Value = 2739.8
$ws.Range("I62").Value = 2739.8
$ws.Range("K62").Value = 2739.8
$ws.Range("M62").Value = -2115.8
$ws.Range("H65").Value = 2739.8
$ws.Range("I65").Value = 2739.8
$ws.Range("K65").Value = 13699
$ws.Range("M65").Value = -10579
$ws.Range("H68").Value = 29673.75
$ws.Range("J68").Value = 29673.75
$ws.Range("L68").Value = 29673.75
$ws.Range("N68").Value = -31171.75
$ws.Range("H71").Value = 29673.75
$ws.Range("J71").Value = 29673.75
$ws.Range("L71").Value = 89021.25
$ws.Range("N71").Value = -96509.25
$ws.Range("H74").Value = 19000
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 33000
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 33000
$ws.Range("M74").Value = -4126
$ws.Range("N74").Value = -34748
$ws.Range("H77").Value = 19000
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 33000
$ws.Range("K77").Value = 15000
$ws.Range("L77").Value = 99000
$ws.Range("M77").Value = -10632
$ws.Range("N77").Value = -107736
$ws.Range("H109").Value = 18940.25
$ws.Range("J109").Value = 18940.25
$ws.Range("L109").Value = 18940.25
$ws.Range("N109").Value = -21020.25
$ws.Range("H134").Value = 4449108.5
$ws.Range("I134").Value = 5163.077
$ws.Range("J134").Value = 33334754
$ws.Range("K134").Value = 15489.231
$ws.Range("L134").Value = 100004262
$ws.Range("M134").Value = -12954.231
$ws.Range("N134").Value = -100009332

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1568.6063
$ws.Range("I68").Value = 1119.2559
$ws.Range("J68").Value = 1947.4706
$ws.Range("K68").Value = 3357.7677
$ws.Range("L68").Value = 5842.4118
$ws.Range("M68").Value = -2546.7677
$ws.Range("N68").Value = -7464.4118
$ws.Range("H71").Value = 1568.6063
$ws.Range("I71").Value = 1119.2559
$ws.Range("J71").Value = 1947.4706
$ws.Range("K71").Value = 10073.3031
$ws.Range("L71").Value = 17527.2354
$ws.Range("M71").Value = -6017.303100000001
$ws.Range("N71").Value = -25639.2354
$ws.Range("H107").Value = 1477.3043
$ws.Range("I107").Value = 922
$ws.Range("J107").Value = 1834.2858
$ws.Range("K107").Value = 2766
$ws.Range("L107").Value = 5502.857400000001
$ws.Range("M107").Value = -846
$ws.Range("N107").Value = -9342.857400000001
$ws.Range("H121").Value = 635850.7
$ws.Range("I121").Value = 366
$ws.Range("J121").Value = 834439.7
$ws.Range("K121").Value = 1098
$ws.Range("L121").Value = 2503319.1
$ws.Range("M121").Value = 212
$ws.Range("N121").Value = -2505939.1

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 24593.666
$ws.Range("J57").Value = 24593.666
$ws.Range("L57").Value = 24593.666
$ws.Range("N57").Value = -26233.666
$ws.Range("H107").Value = 960.619
$ws.Range("I107").Value = 791.5714
$ws.Range("J107").Value = 1298.7142
$ws.Range("K107").Value = 791.5714
$ws.Range("L107").Value = 1298.7142
$ws.Range("M107").Value = 1128.4286
$ws.Range("N107").Value = -5138.7142
$ws.Range("H122").Value = 6475.421
$ws.Range("I122").Value = 9101.75
$ws.Range("J122").Value = 1973.1428
$ws.Range("K122").Value = 27305.25
$ws.Range("L122").Value = 5919.428400000001
$ws.Range("M122").Value = -24855.25
$ws.Range("N122").Value = -10819.4284
$ws.Range("H126").Value = 12043.852
$ws.Range("I126").Value = 2297.8572
$ws.Range("J126").Value = 22539.54
$ws.Range("K126").Value = 6893.571599999999
$ws.Range("L126").Value = 67618.62
$ws.Range("M126").Value = -4423.571599999999
$ws.Range("N126").Value = -72558.62
$ws.Range("H132").Value = 2787.2974
$ws.Range("I132").Value = 2085.9
$ws.Range("J132").Value = 3612.4707
$ws.Range("K132").Value = 6257.700000000001
$ws.Range("L132").Value = 10837.4121
$ws.Range("M132").Value = -3727.700000000001
$ws.Range("N132").Value = -15897.4121

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1367.4667
$ws.Range("I46").Value = 1475
$ws.Range("J46").Value = 1244.5714
$ws.Range("K46").Value = 1475
$ws.Range("L46").Value = 1244.5714
$ws.Range("M46").Value = -1287
$ws.Range("N46").Value = -1620.5714
$ws.Range("H68").Value = 2860.4285
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 3204.6
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 3204.6
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -4702.6
$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26622
$ws.Range("H71").Value = 2860.4285
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 3204.6
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 16023
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -23511
$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -83112
$ws.Range("H136").Value = 1199.2858
$ws.Range("I136").Value = 1211.3334
$ws.Range("J136").Value = 1169.1666
$ws.Range("K136").Value = 3634.0002
$ws.Range("L136").Value = 3507.4998
$ws.Range("M136").Value = -1084.0002
$ws.Range("N136").Value = -8607.4998

